# AWS Config Organization Rules Architecture diagram refresh (Issue #28)
#
# 1) Re-cache the auto "today" date placeholders (slide master, every
#    custom layout, and the notes master) from 9/10/20 -> 5/27/21.
# 2) Drop the translucent slate-blue overlay fill (5A6B86 @ ~9.8% alpha)
#    on the four background "region" rectangles, leaving them unfilled
#    (their outline stays untouched).

$p = $ppt.ActivePresentation
$newDate = "5/27/21"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every custom (slide) layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# Notes master's date placeholder (its Shapes collection is addressed
# through the headers/footers date-and-time field here).
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# Remove the overlay fill from the four region background rectangles.
$slide = $p.Slides.Item(1)
$overlayRectangles = @("Rectangle 128", "Rectangle 126", "Rectangle 127", "Rectangle 57")
foreach ($rectName in $overlayRectangles) {
    $rect = $slide.Shapes.Item($rectName)
    $rect.Fill.Visible = $false
}
